# The "Förändrad" (changed) date column C for rows 2-44 was bumped by
# one day (serial date 45188 -> 45189, i.e. 2023-09-19 -> 2023-09-20).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C44").Value = 45189
